# Sample Project / Main.xlsx — "Rules" sheet edit
#
# Rule row 11 (the "R40" rule) has its Rule-name cell (B11) renamed from
# the text "R40" to the text "1". The new label is stored as a shared
# string (text), not a number, so the leading single-quote forces Excel
# to keep it as text instead of inferring a numeric value.

$wb = $excel.ActiveWorkbook

try {
    $ws = $wb.Worksheets.Item("Rules")
} catch {
    $ws = $wb.ActiveSheet
}

$ws.Range("B11").Value = "'1"
